# Multiply actual emission change in with 100 in backend instead of in
# frontend, to be consistent with needed emission change in percent.
# Column V = "actualEmissionChangePercent"; data rows 2..291.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 22).End(-4162).Row
$colV = 22

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, $colV)
    $val = $cell.Value()
    if ($val -ne $null) {
        $cell.Value() = $val * 100
    }
}
